# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to the Leve profit calculation sheets
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ columns).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15175
$ws.Range("I18").Value = 350
$ws.Range("K18").Value = 350
$ws.Range("M18").Value = -66
$ws.Range("H113").Value = 35501.832
$ws.Range("J113").Value = 3003
$ws.Range("L113").Value = 3003
$ws.Range("N113").Value = -9511
$ws.Range("H129").Value = 2003.0667
$ws.Range("J129").Value = 2587.3635
$ws.Range("L129").Value = 7762.0905
$ws.Range("N129").Value = -17762.0905
$ws.Range("H131").Value = 1352.15
$ws.Range("I131").Value = 947.9167
$ws.Range("J131").Value = 1958.5
$ws.Range("K131").Value = 2843.7501
$ws.Range("L131").Value = 5875.5
$ws.Range("M131").Value = 2196.2499
$ws.Range("N131").Value = -15955.5
$ws.Range("H138").Value = 4782.05
$ws.Range("I138").Value = 5863.2856
$ws.Range("J138").Value = 4552.697
$ws.Range("K138").Value = 17589.8568
$ws.Range("L138").Value = 13658.091
$ws.Range("M138").Value = -12449.8568
$ws.Range("N138").Value = -23938.091

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10467.2
$ws.Range("I32").Value = 9191.071
$ws.Range("K32").Value = 9191.071
$ws.Range("M32").Value = -8904.071
$ws.Range("H139").Value = 45315.8
$ws.Range("J139").Value = 45315.8
$ws.Range("L139").Value = 45315.8
$ws.Range("N139").Value = -55595.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2509.8975
$ws.Range("J20").Value = 2884.8
$ws.Range("L20").Value = 2884.8
$ws.Range("N20").Value = -3378.8
$ws.Range("H58").Value = 19066
$ws.Range("J58").Value = 19066
$ws.Range("L58").Value = 19066
$ws.Range("N58").Value = -19654
$ws.Range("H82").Value = 12611.667
$ws.Range("J82").Value = 48283
$ws.Range("L82").Value = 48283
$ws.Range("N82").Value = -49049
$ws.Range("H85").Value = 12611.667
$ws.Range("J85").Value = 48283
$ws.Range("L85").Value = 48283
$ws.Range("N85").Value = -50935
$ws.Range("H86").Value = 223388.44
$ws.Range("I86").Value = 1333
$ws.Range("J86").Value = 667499.3
$ws.Range("K86").Value = 1333
$ws.Range("L86").Value = 667499.3
$ws.Range("M86").Value = -210
$ws.Range("N86").Value = -669745.3
$ws.Range("H89").Value = 223388.44
$ws.Range("I89").Value = 1333
$ws.Range("J89").Value = 667499.3
$ws.Range("K89").Value = 6665
$ws.Range("L89").Value = 3337496.5
$ws.Range("M89").Value = -1049
$ws.Range("N89").Value = -3348728.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1978230.5
$ws.Range("I58").Value = 2719051.8
$ws.Range("J58").Value = 2706.8333
$ws.Range("K58").Value = 2719051.8
$ws.Range("L58").Value = 2706.8333
$ws.Range("M58").Value = -2718848.8
$ws.Range("N58").Value = -3112.8333
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H134").Value = 1218.0256
$ws.Range("I134").Value = 1079.6333
$ws.Range("J134").Value = 1679.3334
$ws.Range("K134").Value = 3238.8999
$ws.Range("L134").Value = 5038.0002
$ws.Range("M134").Value = -703.8998999999999
$ws.Range("N134").Value = -10108.0002
$ws.Range("H136").Value = 1978230.5
$ws.Range("I136").Value = 2719051.8
$ws.Range("J136").Value = 2706.8333
$ws.Range("K136").Value = 8157155.399999999
$ws.Range("L136").Value = 8120.499899999999
$ws.Range("M136").Value = -8154605.399999999
$ws.Range("N136").Value = -13220.4999
$ws.Range("M99").ClearContents()
$ws.Range("M126").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5500
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 3000
$ws.Range("M3").Value = -2888
$ws.Range("H61").Value = 483.33334
$ws.Range("I61").Value = 483.33334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1450.00002
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1235.00002
$ws.Range("H97").Value = 988.75
$ws.Range("I97").Value = 235
$ws.Range("K97").Value = 705
$ws.Range("M97").Value = -209
$ws.Range("H98").Value = 1835.7273
$ws.Range("I98").Value = 699.75
$ws.Range("J98").Value = 2484.8572
$ws.Range("K98").Value = 2099.25
$ws.Range("L98").Value = 7454.571599999999
$ws.Range("M98").Value = -601.25
$ws.Range("N98").Value = -10450.5716
$ws.Range("H115").Value = 4877.9287
$ws.Range("I115").Value = 2249.75
$ws.Range("K115").Value = 6749.25
$ws.Range("M115").Value = -5574.25
$ws.Range("H122").Value = 1949.3334
$ws.Range("J122").Value = 1949.3334
$ws.Range("L122").Value = 17544.0006
$ws.Range("N122").Value = -22444.0006
$ws.Range("H131").Value = 11544.838
$ws.Range("J131").Value = 11544.838
$ws.Range("L131").Value = 34634.514
$ws.Range("N131").Value = -44714.514
$ws.Range("H137").Value = 7219.696
$ws.Range("J137").Value = 8524.352999999999
$ws.Range("L137").Value = 25573.059
$ws.Range("N137").Value = -35773.05899999999
$ws.Range("N61").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2025.2
$ws.Range("I122").Value = 2042.3334
$ws.Range("K122").Value = 6127.0002
$ws.Range("M122").Value = -3677.0002
$ws.Range("H135").Value = 79920
$ws.Range("J135").Value = 79920
$ws.Range("L135").Value = 79920
$ws.Range("N135").Value = -90060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 3800
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3800
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3800
$ws.Range("N17").Value = -4140
$ws.Range("H21").Value = 10000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10348
$ws.Range("H22").Value = 1450.2
$ws.Range("I22").Value = 957.8333
$ws.Range("K22").Value = 957.8333
$ws.Range("M22").Value = -662.8333
$ws.Range("H27").Value = 1450.2
$ws.Range("I27").Value = 957.8333
$ws.Range("K27").Value = 957.8333
$ws.Range("M27").Value = -850.8333
$ws.Range("H55").Value = 5263528.5
$ws.Range("J55").Value = 411.77777
$ws.Range("L55").Value = 411.77777
$ws.Range("N55").Value = -757.7777699999999
$ws.Range("M17").ClearContents()
$ws.Range("M21").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 40403.332
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 40403.332
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 40403.332
$ws.Range("N7").Value = -40629.332
$ws.Range("H15").Value = 9987.5
$ws.Range("J15").Value = 9987.5
$ws.Range("L15").Value = 9987.5
$ws.Range("N15").Value = -10563.5
$ws.Range("H19").Value = 10613
$ws.Range("J19").Value = 11866.667
$ws.Range("L19").Value = 11866.667
$ws.Range("N19").Value = -12214.667
$ws.Range("H63").Value = 50000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 50000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51248
$ws.Range("H66").Value = 50000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 50000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156240
$ws.Range("H132").Value = 1615.625
$ws.Range("I132").Value = 988.0526
$ws.Range("J132").Value = 4000.4
$ws.Range("K132").Value = 2964.1578
$ws.Range("L132").Value = 12001.2
$ws.Range("M132").Value = -434.1578
$ws.Range("N132").Value = -17061.2
$ws.Range("H136").Value = 37041130
$ws.Range("I136").Value = 61731324
$ws.Range("K136").Value = 185193972
$ws.Range("M136").Value = -185191422
$ws.Range("M7").ClearContents()
$ws.Range("M63").ClearContents()
$ws.Range("M66").ClearContents()
